$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.723.66"
$ws.Range("D3").Value = "1.540.32"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'290.08"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").Value = "'0.3903"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "'0.3174"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "'42.89"
$ws.Range("E9").Value = "  +4.44%  "
$ws.Range("D10").Value = "'0.07184"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'1.056"
$ws.Range("E11").Value = "  -5.87%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'5.626"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.620"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.557.78"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "'0.06577"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "'83.16"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "'6.151"
$ws.Range("E21").Value = "  -4.20%  "
$ws.Range("D22").Value = "'15.37"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("E23").Value = "  -5.60%  "
$ws.Range("E24").Value = "  +7.23%  "
$ws.Range("D25").Value = "21.729.29"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'2.360"
$ws.Range("E26").Value = "  -6.24%  "
$ws.Range("D27").Value = "'146.51"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "'18.37"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "'4.842"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "1.718.98"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'117.49"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("D32").Value = "'5.895"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'0.9642"
$ws.Range("E33").Value = "  -13.89%  "
$ws.Range("D34").Value = "'0.08192"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "'8.803"
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").Value = "'0.06082"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "'5.120"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").Value = "'0.02198"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").Value = "'0.2034"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.179"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.430"
$ws.Range("E41").Value = "  -12.52%  "
$ws.Range("D43").Value = "'10.64"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'0.5724"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.740"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.98"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").Value = "'0.5484"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("D48").Value = "'1.158"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").Value = "'116.33"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Value = "'1.869"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("D51").Value = "'0.06710"
$ws.Range("E51").Value = "  -2.77%  "
